# "Generate Report for Handoff" - refresh the localization-status report:
# flip the in-progress rows to "Ready for handoff" and bump their
# handoff/generate timestamps, then resize the now-wider "Status"/status
# columns to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: per-locale status + the latest handoff-xliff generate date
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-16 04:55:49"

# zh-cn detail sheet: Status + Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-16 04:55:45"

# de-de detail sheet: Status + Latest Handoff Datetime
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-16 04:55:49"

# Widen the columns that now hold the longer "Ready for handoff" text.
$overview.Range("E1").EntireColumn.ColumnWidth = 16.333333333333332
$overview.Range("F1").EntireColumn.ColumnWidth = 16.333333333333332
$zhcn.Range("C1").EntireColumn.ColumnWidth = 16.333333333333332
$dede.Range("C1").EntireColumn.ColumnWidth = 16.333333333333332
